$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The figures in this sheet are stored as plain text (not numbers), e.g. "1,035".
# Setting .Value directly on a string containing a comma makes the engine treat it
# as a formatted number (applying a "#,##0" style) instead of literal text, so we
# force text interpretation via NumberFormat "@" and then restore the cell style
# back to Normal/General so no stray numeric formatting is left behind.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 11 - Retained Earnings
Set-TextValue $ws.Range("C11") "1,055"
Set-TextValue $ws.Range("D11") "2,200"
Set-TextValue $ws.Range("E11") "3,553"
Set-TextValue $ws.Range("F11") "5,154"
Set-TextValue $ws.Range("G11") "7,064"

# Row 12 - Total Equity
Set-TextValue $ws.Range("C12") "1,155"
Set-TextValue $ws.Range("D12") "2,300"
Set-TextValue $ws.Range("E12") "3,653"
Set-TextValue $ws.Range("F12") "5,254"
Set-TextValue $ws.Range("G12") "7,164"

# Row 13 - Total Liabilities & Equity
Set-TextValue $ws.Range("C13") "1,385"
Set-TextValue $ws.Range("D13") "2,530"
Set-TextValue $ws.Range("E13") "3,892"
Set-TextValue $ws.Range("F13") "5,503"
Set-TextValue $ws.Range("G13") "7,424"
